$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid_Login")

# Row 6: flip the RunToTest flag from Y to N (logistic type changed)
$ws.Range("A6").Value = "N"

# Refresh the existing Surekha password value (shared string backing C6 / its hyperlink)
$ws.Range("C6").Value = "Surekha@003"

# Add a new row (row 8) of login test data for Nasreen
$ws.Range("A8").Value = "Y"
$ws.Range("B8").Value = "Nasreen.khan@ram.co.za"
$ws.Range("C8").Value = "Nasreen@1235"

$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:Nasreen.khan@ram.co.za")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:Nasreen@1235")

$ws.Range("B8").Style = "Hyperlink"
$ws.Range("C8").Style = "Hyperlink"

$ws.Range("B10").Select()
